# Importer module sheet: replace the hard-coded "Nom/Email/Classe" sample
# rows with the NomModule/Ensiegnant_Email/Classes data used by the
# import-module modal, and refresh the professor e-mail hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -----------------------------------------------------
$ws.Range("D2").Value = "NomModule"
$ws.Range("E2").Value = "Ensiegnant_Email"
$ws.Range("F2").Value = "Classes"

# --- data rows --------------------------------------------------------
$ws.Range("D3").Value = "Java"
$ws.Range("F3").Value = "3.GInfo"

$ws.Range("D4").Value = "T.I"
$ws.Range("F4").Value = "4.Ginfo/4.GTR"

$ws.Range("D5").Value = "R.O"
$ws.Range("F5").Value = "4.Ginfo/4.GTR"

$ws.Range("D6").Value = "C#"
$ws.Range("F6").Value = "4.Ginfo"

$ws.Range("D7").Value = "Analyse 3"
$ws.Range("F7").Value = "Cp 2"

# --- teacher e-mails + hyperlinks --------------------------------------
# Drop every existing hyperlink (and its relationship) first so the
# rebuilt set gets fresh rIds in the right order.
$ws.Cells.Hyperlinks.Delete()

# E3 keeps the "Lien hypertexte" look but is no longer a live link.
$ws.Range("E3").Value = "profmail1@gmail.com"
$ws.Range("E3").Style = "Lien hypertexte"

$ws.Range("E4").Value = "profmail2@gmail.com"
$ws.Range("E5").Value = "profmail3@gmail.com"
$ws.Range("E6").Value = "profmail4@gmail.com"
$ws.Range("E7").Value = "profmail5@gmail.com"

$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:profmail2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:profmail3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:profmail4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:profmail5@gmail.com")

# Re-apply the built-in hyperlink style so Add()'s extra style variant
# collapses back onto the workbook's existing "Lien hypertexte" xf.
$ws.Range("E4:E7").Style = "Lien hypertexte"

# --- column widths (author widened D:F and dropped the old bestFit) ----
$ws.Columns("D").ColumnWidth = 13.333333333333334
$ws.Columns("E").ColumnWidth = 22
$ws.Columns("F").ColumnWidth = 18.833333333333332

# --- selection moved from K9 to H9 -------------------------------------
$ws.Range("H9").Select()
